# Update the "Salida" column (column B, rows 2-81) on the "Solucion" sheet
# with the new assignment order produced by the beam-search run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Solucion")

$salidas = @(
    "S052","S072","S040","S037","S030","S054","S057","S018","S022","S055",
    "S028","S068","S069","S015","S032","S001","S011","S025","S005","S034",
    "S026","S076","S067","S016","S021","S035","S047","S033","S004","S023",
    "S079","S006","S075","S074","S017","S073","S031","S002","S019","S064",
    "S080","S048","S071","S077","S066","S039","S038","S053","S014","S060",
    "S050","S063","S042","S049","S061","S058","S020","S056","S029","S051",
    "S013","S043","S009","S065","S041","S045","S036","S062","S010","S059",
    "S070","S003","S008","S078","S024","S044","S046","S007","S012","S027"
)

for ($i = 0; $i -lt $salidas.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $salidas[$i]
}
